# Update "peringkat" data: refreshed ranking rows (ID, Kualitas Servis,
# Harga, Skor Kelayakan) after refactoring inferensi/defuzzifikasi rules.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  79, 92, 22360, 90),
    @(3,  80, 89, 22012, 89),
    @(4,  78, 86, 27315, 86),
    @(5,  69, 85, 24551, 85),
    @(6,  86, 84, 29811, 84),
    @(7,  25, 94, 34513, 80.97399999999999),
    @(8,  66, 80, 20052, 80),
    @(9,  88, 100, 35304, 79.392),
    @(10, 23, 77, 22825, 77),
    @(11, 33, 73, 24704, 73)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
